# Fix situation calcul issues:
# - Insert two new rows (new contracts "512/CASA 2") above the existing data
# - Rename the "903/CASA ANFA" contract to "903/CASA ANFA/AV"
# - Correct a couple of tax amounts that were miscalculated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2-10) down by two rows to make room
# for the two new rows at the top of the table.
$ws.Range("A2:M3").EntireRow.Insert()

# New row 2: "512/CASA 2" / STE LOCATION
$ws.Range("A2").Value = "512/CASA 2"
$ws.Range("B2").Value = "Point de vente"
$ws.Range("C2").Value = "31451"
$ws.Range("D2").Value = "STE LOCATION "
$ws.Range("E2").Value = "ds"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 5000
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 5000

# New row 3: "512/CASA 2" / STE MAISON
$ws.Range("A3").Value = "512/CASA 2"
$ws.Range("B3").Value = "Point de vente"
$ws.Range("C3").Value = "56987"
$ws.Range("D3").Value = "STE MAISON "
$ws.Range("E3").Value = "ds"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 5000
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 5000

# Former row 2 (now row 4): rename contract to "903/CASA ANFA/AV"
$ws.Range("A4").Value = "903/CASA ANFA/AV"

# Former row 8 (now row 10): fix the Taxe/loyer amount (10% of 28000 = 2800, not 8400)
$ws.Range("J10").Value = 2800

# Former row 10 (now row 12): fix the Taxe/loyer amount (15% of 150000 = 22500, not 270000)
$ws.Range("J12").Value = 22500
